$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column M ("Issue_Cause_Division_Act") mirrors the formatting of the
# existing "Issue_Cause_Division" column (L), cell by cell, then gets the
# actual division value (collapsing the combined "... and Campaign Dev"
# values down to the owning division).
$values = @{
    1  = "Issue_Cause_Division_Act"
    2  = "COPS"
    3  = "COPS"
    4  = "COPS"
    5  = "COPS"
    6  = "COPS"
    7  = "COPS"
    8  = "COPS"
    9  = "COPS"
    10 = "COPS"
    11 = "EMOPs"
    12 = "EMOPs"
    13 = "DMOPs"
}

for ($row = 1; $row -le 13; $row++) {
    $ws.Cells.Item($row, 12).Copy()
    $ws.Cells.Item($row, 13).PasteSpecial(-4122)
    $ws.Cells.Item($row, 13).Value = $values[$row]
}

$excel.CutCopyMode = 0

# Size the new column to fit its contents, like the other columns.
$ws.Columns("M").AutoFit()

# Update the active selection to reflect where the user left off editing.
$ws.Range("K7").Select()
